$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting (and placeholder values) from row 30 into the two new rows.
$ws.Range("A30:E30").Copy($ws.Range("A31:E31"))
$ws.Range("A30:E30").Copy($ws.Range("A32:E32"))

# Row 31 keeps the same taller row height as row 30 (28.8pt); row 32 uses the
# default sheet row height.
$ws.Rows.Item(31).RowHeight = 28.8
$ws.Rows.Item(32).UseStandardHeight = $true

# Fill in the new cell values in the same order the original author typed them,
# so that new shared-string entries are appended in the expected order.
$ws.Range("A31").Value = "IAM030"
$ws.Range("A32").Value = "IAM031"
$ws.Range("B32").Value = "OPQA-2838"
$ws.Range("C32").Value = "Verify that deep linking is working correctly for help page using FB and LI accounts"
$ws.Range("B31").Value = "OPQA-2837"
$ws.Range("C31").Value = "Verify that deep linking is working correctly for account page using FB and LI accounts"
$ws.Range("D31").Value = "Y"
$ws.Range("E31").Value = "PASS"
$ws.Range("D32").Value = "Y"
$ws.Range("E32").Value = "PASS"

# Update the selection to match the final state (active cell C32).
$ws.Range("C32").Select()
